$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '30.501.43'
$ws.Cells.Item(2, 5).Value2 = '  -0.07%  '

$ws.Cells.Item(3, 4).Value2 = '1.872.50'
$ws.Cells.Item(3, 5).Value2 = '  -0.81%  '

$ws.Cells.Item(4, 5).Value2 = '  +0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '247.61'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value2 = '  +1.08%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = '1.000'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value2 = '  +0.13%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = '0.4734'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value2 = '  -0.73%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = '0.2895'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value2 = '  -0.10%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = '0.06464'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value2 = '  -1.17%  '

$ws.Cells.Item(10, 5).Value2 = '  +3.20%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '0.07693'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value2 = '  -1.03%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = '0.7385'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value2 = '  -0.04%  '

$ws.Cells.Item(13, 2).Value2 = 'Litecoin'
$ws.Cells.Item(13, 3).Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = '95.99'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value2 = '  -0.29%  '

$ws.Cells.Item(14, 2).Value2 = 'WrappedEther'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value2 = '1.871.27'
$ws.Cells.Item(14, 5).Value2 = '  -0.90%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = '5.168'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value2 = '  +0.03%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '274.57'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value2 = '  -0.61%  '

$ws.Cells.Item(17, 4).Value2 = '30.535.72'
$ws.Cells.Item(17, 5).Value2 = '  +0.08%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = '13.24'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value2 = '  -2.10%  '

$ws.Cells.Item(19, 5).Value2 = '  +0.07%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = '0.000007464'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value2 = '  -1.75%  '

$ws.Cells.Item(21, 4).Value2 = '2.108.26'
$ws.Cells.Item(21, 5).Value2 = '  -1.25%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '1.001'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value2 = '  +0.12%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = '5.231'
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value2 = '  -1.33%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value2 = '6.162'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value2 = '  -0.75%  '

$ws.Cells.Item(25, 2).Value2 = 'Monero'
$ws.Cells.Item(25, 3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '165.03'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value2 = '  -0.08%  '

$ws.Cells.Item(26, 2).Value2 = 'Cosmos'
$ws.Cells.Item(26, 3).Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = '9.178'
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value2 = '  -1.43%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = '18.66'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value2 = '  -2.00%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = '1.898'
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value2 = '  -3.74%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = '0.09954'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value2 = '  -0.13%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = '1.343'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value2 = '  -3.12%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = '1.507'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value2 = '  -0.46%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = '4.226'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value2 = '  -2.94%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = '4.077'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value2 = '  -0.33%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = '0.04758'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value2 = '  -0.53%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = '1.116'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value2 = '  -1.56%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = '0.6909'
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value2 = '  -1.69%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = '2.716'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value2 = '  -0.09%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value2 = '0.01854'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value2 = '  -0.11%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '2.754'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value2 = '  -0.42%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '6.233'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value2 = '  -3.83%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = '72.97'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value2 = '  +3.21%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = '1.962'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value2 = '  +1.76%  '

$ws.Cells.Item(43, 5).Value2 = '  +0.11%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = '0.4146'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value2 = '  -0.92%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = '0.8332'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value2 = '  -1.48%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '101.33'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value2 = '  -1.49%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = '9.332'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value2 = '  -0.78%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = '35.26'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value2 = '  -0.22%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = '6.958'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value2 = '  -2.90%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = '914.11'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value2 = '  -2.14%  '

$ws.Cells.Item(51, 5).Value2 = '  +1.02%  '
